$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.345.78'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '2.370.94'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.701'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.65'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.78%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +26.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.103'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '31.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +16.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +18.87%  '
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = '2.724.33'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '17.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.921'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.36%  '
$ws.Range('D17').Value = '2.370.45'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '44.414.43'
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '79.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '259.21'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.02%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.59%  '
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('E32').Value = '  +7.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.39'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0762'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.47'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.91%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  +8.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.33'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.55%  '
$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.200'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +18.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.102'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.30%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.49%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.87%  '
$ws.Range('E47').Value = '  +11.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = '1.470.20'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.74%  '
